# Daily attendance processing - 2025-11-05 15:47:03
# Normalizes the "Recorded By" (column G) lists so the "System" entry
# (exact case) always sorts to the end of the comma-separated list,
# leaving every other entry (including a differently-cased "system")
# in its original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Test-CaseSensitiveEquals($a, $b) {
    if ($a.Length -ne $b.Length) { return $false }
    for ($i = 0; $i -lt $a.Length; $i++) {
        if ([int][char]$a[$i] -ne [int][char]$b[$i]) {
            return $false
        }
    }
    return $true
}

function Move-SystemToEnd($s) {
    $parts = @($s -split ", ")
    $newParts = @()
    $foundSystem = $false
    foreach ($p in $parts) {
        if (Test-CaseSensitiveEquals $p "System") {
            $foundSystem = $true
        } else {
            $newParts += $p
        }
    }
    if ($foundSystem) {
        $newParts += "System"
    }
    return ($newParts -join ", ")
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column G holds "Recorded By" per the header row; data starts on row 2.
$col = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Value2
    if ([string]::IsNullOrEmpty($current)) {
        continue
    }
    if ($current.IndexOf("System") -lt 0) {
        continue
    }
    $updated = Move-SystemToEnd $current
    if (-not (Test-CaseSensitiveEquals $updated $current)) {
        $cell.Value2 = $updated
    }
}
